$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 11
$ws.Range("B113").Value = "Vega Monumental Concepción"
$ws.Range("C113").Value = "Bíobío"
$ws.Range("D113").Value = 45072
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = 100112043
$ws.Range("G113").Value = "Pepino ensalada"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 220
$ws.Range("K113").Value = 12000
$ws.Range("L113").Value = 14000
$ws.Range("M113").Value = 12909
$ws.Range("N113").Value = "$/caja 60 unidades"
$ws.Range("O113").Value = "Región de Arica y Parinacota"
$ws.Range("P113").Value = 215
$ws.Range("Q113").Value = 60
$ws.Range("R113").Value = "Hortaliza"
